# "brute force algorithm initialized"
#
# Tanarok (Teachers) sheet: the first two teachers (Kerekes/Janos and
# Tolmacsi/Agnes) are removed, leaving only Nagy/Jozsef and Kis/Balint
# (previously rows 3-4, now rows 1-2). Their time-slot list text is
# normalized to drop the spaces after the commas.
#
# Diakok (Students) sheet: every student's "C" column (which used to
# hold the full list of compatible time slots) is collapsed down to a
# single, concretely-assigned time slot - the result of a first brute
# force matching pass.
#
# Selections/active sheet are updated to match where the author's
# cursor ended up after these edits.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tanarok")
$ws2 = $wb.Worksheets.Item("Diakok")

# --- Tanarok: drop the first two teachers, keep the other two ---
[void]$ws1.Rows.Item(1).Delete()
[void]$ws1.Rows.Item(1).Delete()

# Remaining rows keep their names, only the schedule text loses its spaces
$ws1.Range("C1").Value = "K5,SZ7,CS2,P1"
$ws1.Range("C2").Value = "H1,H2,H3,K1,K2,K3"

# --- Diakok: assign each student a single concrete time slot ---
$ws2.Range("C1").Value = "K5"
$ws2.Range("C2").Value = "SZ7"
$ws2.Range("C3").Value = "H1"
$ws2.Range("C4").Value = "H2"

# --- Restore the view/selection state from the edited workbook ---
[void]$ws2.Range("C12").Select()
[void]$ws1.Range("C6").Select()
[void]$ws1.Activate()
